$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.369.79"
$ws.Range("E2").Value = "  -1.14%  "

$ws.Range("D3").Value = "3.247.89"
$ws.Range("E3").Value = "  +3.04%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.31%  "

$ws.Range("D8").Value = "3.249.07"
$ws.Range("E8").Value = "  +3.21%  "

$ws.Range("E9").Value = "  -1.80%  "

$ws.Range("E10").Value = "  -1.32%  "

$ws.Range("E11").Value = "  -0.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.465"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.64%  "

$ws.Range("E13").Value = "  -3.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.73%  "

$ws.Range("D15").Value = "3.777.97"
$ws.Range("E15").Value = "  +2.98%  "

$ws.Range("E16").Value = "  -0.39%  "

$ws.Range("D17").Value = "3.243.68"
$ws.Range("E17").Value = "  +3.33%  "

$ws.Range("D18").Value = "63.367.48"
$ws.Range("E18").Value = "  -1.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "474.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.72%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.31%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.733"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.95%  "

$ws.Range("E25").Value = "  -0.59%  "

$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.80%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.12"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.58%  "

$ws.Range("E30").Value = "  +2.65%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.53"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.11%  "

$ws.Range("E32").Value = "  -0.17%  "

$ws.Range("E33").Value = "  -4.63%  "

$ws.Range("E34").Value = "  -4.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.24%  "

$ws.Range("E37").Value = "  -0.18%  "

$ws.Range("E38").Value = "  -4.70%  "

$ws.Range("E39").Value = "  -1.23%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "422.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.53%  "

$ws.Range("E41").Value = "  +0.06%  "

$ws.Range("D42").Value = "2.980.89"
$ws.Range("E42").Value = "  +2.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.54%  "

$ws.Range("E44").Value = "  -8.18%  "

$ws.Range("E45").Value = "  +1.86%  "

$ws.Range("E46").Value = "  -0.93%  "

$ws.Range("E47").Value = "  +0.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.15%  "

$ws.Range("E49").Value = "  -3.66%  "

$ws.Range("E50").Value = "  -0.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.44%  "
